{"js": "// Update the worksheet date heading and the 25 division-problem answers\n// in the 5x5 table, preserving all existing run/paragraph formatting.\n// Each new value replaces the old one strictly by position (row/column),\n// since several old answers are not unique strings within the document.\n\nconst body = context.document.body;\n\n// 1) Date heading paragraph (first paragraph of the document, outside the table).\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst headingRange = paras.items[0].getRange();\nheadingRange.insertText(\"2023-11-22 Wednesday\", Word.InsertLocation.replace);\n\n// 2) Table with the division problems.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// New values for each of the 5 \"content\" rows (table rows 0, 4, 8, 12, 16 \u2014\n// the rows in between are intentionally blank spacer rows).\nconst newValues = [\n  [\"59\u00f73=19, 2\", \"51\u00f74=12, 3\", \"47\u00f79=5, 2\", \"15\u00f72=7, 1\", \"79\u00f76=13, 1\"],\n  [\"38\u00f77=5, 3\", \"74\u00f74=18, 2\", \"80\u00f74=20, 0\", \"33\u00f72=16, 1\", \"68\u00f75=13, 3\"],\n  [\"40\u00f75=8, 0\", \"61\u00f77=8, 5\", \"47\u00f76=7, 5\", \"70\u00f75=14, 0\", \"18\u00f76=3, 0\"],\n  [\"92\u00f78=11, 4\", \"80\u00f72=40, 0\", \"70\u00f79=7, 7\", \"16\u00f73=5, 1\", \"89\u00f79=9, 8\"],\n  [\"91\u00f74=22, 3\", \"83\u00f79=9, 2\", \"81\u00f74=20, 1\", \"70\u00f77=10, 0\", \"17\u00f78=2, 1\"],\n];\nconst contentRowIndexes = [0, 4, 8, 12, 16];\n\nfor (let r = 0; r < contentRowIndexes.length; r++) {\n  const tableRowIndex = contentRowIndexes[r];\n  for (let c = 0; c < 5; c++) {\n    const cell = table.getCell(tableRowIndex, c);\n    const cellPara = cell.body.paragraphs.getFirst();\n    const cellRange = cellPara.getRange();\n    cellRange.insertText(newValues[r][c], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date heading and the 25 division-problem answers\n# in the 5x5 table, preserving all existing run/paragraph formatting.\n# Each new value replaces the old one strictly by position (row/column),\n# since several old answers are not unique strings within the document.\n\n$d = $word.ActiveDocument\n\n# 1) Date heading paragraph (first paragraph of the document, outside the table).\n$d.Paragraphs.Item(1).Range.Text = \"2023-11-22 Wednesday\"\n\n# 2) Table with the division problems.\n$t = $d.Tables.Item(1)\n\n# New values for each of the 5 \"content\" rows (table rows 1, 5, 9, 13, 17 \u2014\n# the rows in between are intentionally blank spacer rows).\n$newValues = @(\n    @(\"59\u00f73=19, 2\", \"51\u00f74=12, 3\", \"47\u00f79=5, 2\", \"15\u00f72=7, 1\", \"79\u00f76=13, 1\"),\n    @(\"38\u00f77=5, 3\", \"74\u00f74=18, 2\", \"80\u00f74=20, 0\", \"33\u00f72=16, 1\", \"68\u00f75=13, 3\"),\n    @(\"40\u00f75=8, 0\", \"61\u00f77=8, 5\", \"47\u00f76=7, 5\", \"70\u00f75=14, 0\", \"18\u00f76=3, 0\"),\n    @(\"92\u00f78=11, 4\", \"80\u00f72=40, 0\", \"70\u00f79=7, 7\", \"16\u00f73=5, 1\", \"89\u00f79=9, 8\"),\n    @(\"91\u00f74=22, 3\", \"83\u00f79=9, 2\", \"81\u00f74=20, 1\", \"70\u00f77=10, 0\", \"17\u00f78=2, 1\")\n)\n$contentRows = @(1, 5, 9, 13, 17)\n\nfor ($r = 0; $r -lt $contentRows.Length; $r++) {\n    $rowIndex = $contentRows[$r]\n    for ($c = 1; $c -le 5; $c++) {\n        $cell = $t.Cell($rowIndex, $c)\n        $cell.Range.Text = $newValues[$r][$c - 1]\n    }\n}\n"}
